# switch to datannur name
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# about_page_1 value (B4) : 'Fonctionnement' section, Datannuaire -> datannur
$fonctionnement = @'
### Fonctionnement
datannur contient 7 entités principales. On peut les diviser en deux catégories, partie intérieur et partie extérieur aux datasets. Le **dataset** représente une table de base de données ou un fichier de données (excel, csv, ...) sous forme de tableau (lignes et colonnes).

mermaid( 
  $dataset -.-> intérieur
  $dataset -.-> extérieur
);

'@

# contact_email value (B2) : dat@nnuaire.email -> contact@datannur.com
$contact = @'
contact@datannur.com
'@

# about_main value (B3) : main banner markdown, datannuaire -> datannur branding
$banner = @'
![main_banner not_rounded](data/img/main_banner{dark_mode}.png?v=1)

# datannur, le catalogue de données portable

Permet de **centraliser**, **rechercher** et **visualiser** les informations sur une collection de jeux de données

Pour améliorer l’organisation des données et faciliter leur **partage** et leur **documentation**

**Simple** et **flexible**, s’intègre rapidement dans tous types d’environnement


- **Facile** :
Aucune installation ou configuration nécessaire, aucun coût ou prérequis technique

- **Portable** :
Fonctionne partout (local, cloud, disque partagé), un simple dossier que l’on peut copier, déplacer, envoyer et ouvrir avec n’importe quel navigateur web

- **Complet** :
Flexible, complet et structuré autour de 6 concepts avec un niveau de détail important : Institution, Dossier, Mot clé, Dataset, Variable et Modalité

- **Indépendant** :
Le catalogue n’est qu’une interface pour visualiser les métadonnées, le processus de leur création et mise à jour est indépendant et sous votre contrôle

- **Sécurisé** :
De pars la séparation stricte entre les deux systèmes, l’application est isolée dans le navigateur, ne peut rien modifier sur la machine et ne pose ainsi aucun risque

La version ici présente est un **prototype** en cours de développement et d'expérimentation. Les données utilisées sont fictives et uniquement à usage de test et de développement. Question ou suggestion : [contact@datannur.com](mailto:contact@datannur.com).
'@

# about_main_meta value (B9) : meta view description, Meta -> meta, link anchor change
$metaDesc = @'
La vue méta représente les méta-méta-données, c’est-à-dire l'information sur les jeux de données internes au catalogue. 

Cette vue apporte une synthèse générale sur le contenu du catalogue et permet de voir l'architecture de données de datannur.

Pour vérifier l'intégriter de la base de données, cliquez sur ce [lien](?app_mode=check_db).
'@

# Apply edits in the same order the shared strings were appended upstream
$ws.Range("B4").Value = $fonctionnement
$ws.Range("B2").Value = $contact
$ws.Range("B3").Value = $banner
$ws.Range("B9").Value = $metaDesc

# Reflect the author's last edit position: active/selected cell moves to B9
$ws.Range("B9").Select()

